$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New student rows to append below the existing header row.
$data = @(
    @(1011825, "DGES", 7, "3er. CICLO EBI", "desconocido", "M", "Estudiante", 20, 0, "", "", "", 12),
    @(1011825, "DGES", 7, "3er. CICLO EBI", "desconocido", "M", "Estudiante", 20, 0, "", "", "", 12),
    @(1011825, "DGES", 8, "3er. CICLO EBI", "desconocido", "M", "Estudiante", 20, 0, "", "", "", 12),
    @(1011825, "DGES", 8, "3er. CICLO EBI", "desconocido", "M", "Estudiante", 20, 0, "", "", "", 12),
    @(1011825, "DGES", 9, "3er. CICLO EBI", "desconocido", "M", "Estudiante", 20, 0, "", "", "", 12),
    @(1011825, "DGES", 9, "3er. CICLO EBI", "desconocido", "M", "Estudiante", 20, 0, "", "", "", 12),
    @(1012812, "DGES", 7, "3er. CICLO EBI", "desconocido", "F", "Estudiante", 18, 0, "", "", "", 94),
    @(1012812, "DGES", 8, "3er. CICLO EBI", "desconocido", "F", "Estudiante", 18, 0, "", "", "", 94),
    @(1012812, "DGES", 9, "3er. CICLO EBI", "desconocido", "F", "Estudiante", 18, 0, "", "", "", 94)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $data[$i]
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $rowValues[$col - 1]
    }
}
